$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values per the diff.
# Cells are stored as text (inlineStr), so we force text assignment
# via NumberFormat "@" before setting .Value to preserve exact formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "305.27"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.24%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.42%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.25%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.49%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.867"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.92%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.133"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-1.14%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.788"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.76%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9218"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.86%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1286"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-7.28%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1902"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.09%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09151"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.25%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03408"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-6.09%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09862"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.61%"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.10%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006163"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "4.41%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.853"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "8.39%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.75%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.39%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.28%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.187"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.85%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-7.93%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.44%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.98%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004886"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.10%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001252"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-19.84%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "42.09%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01945"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.76%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05171"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.70%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007651"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.14%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "10.45%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.61%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002153"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009627"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-15.38%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006174"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.36%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.04%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "65.11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.43%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "39.40%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.04%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.04%"

$wb.Save()
